$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: same look as the other header cells (bold/centered/bordered).
# Copy G1's formatting (style) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data cell H2 with numeric value 0 (plain/default formatting, like the other data cells)
$ws.Range("H2").Value = 0
